# Apply the "Updated symbol list" edits to Sheet1.
# Column D holds numeric-looking price strings that must remain text (inlineStr),
# so we briefly force a text NumberFormat while assigning the value and then restore
# the default "Normal" style so the cell formatting is unchanged from before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "244.26"
Set-TextValue "D3" "25.03"
Set-TextValue "D4" "5.124"
Set-TextValue "D6" "6.465"
Set-TextValue "D7" "3.117"
Set-TextValue "D8" "0.8097"
Set-TextValue "D9" "0.8405"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1338"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.06938"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03130"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "ProBitToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D13" "0.1337"
$ws.Range("E13").Value = "12ProBitTokenPROB"
Set-TextValue "D14" "0.02816"
Set-TextValue "D15" "0.09372"
Set-TextValue "D16" "3.746"
Set-TextValue "D17" "0.001510"
Set-TextValue "D18" "0.04652"
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D19" "0.0005987"
$ws.Range("E19").Value = "18OneONE"
$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D20" "0.006105"
$ws.Range("E20").Value = "19TigerCashTCH"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D21" "0.001238"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "HotbitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D22" "0.004280"
$ws.Range("E22").Value = "21HotbitTokenHTB"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D23" "0.00008696"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D24" "3.502"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D25" "2.109"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D26" "0.3175"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
Set-TextValue "D28" "0.0002317"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1052"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002889"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.006377"
$ws.Range("E43").Value = "42KickTokenKICK"
Set-TextValue "D44" "0.007383"
Set-TextValue "D45" "0.00005302"
Set-TextValue "D47" "0.2519"
Set-TextValue "D48" "0.002275"

Write-Output "Applied 77 cell updates to Sheet1."
